# This script applies the "gth and rgth file" commit:
#  - Renames the "demo" sheet to "rgth" and fills it with a
#    "Receive Gift Auto" request form.
#  - Adds a new "pgth" sheet with a "Provide Gift Auto" request form and
#    makes it the active sheet/tab.
#  - Updates the "dpc" sheet's request data (COI_Auto) and adds a
#    mailto hyperlink.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "dpc" sheet
# ---------------------------------------------------------------------
$dpc = $wb.Worksheets.Item("dpc")

$dpc.Range("B2").Value = "COI_Auto"

$dpc.Range("A3").Value = "competitor"
$dpc.Range("B3").Value = "ak@gmail.com"

$dpc.Range("A4").Value = "company competes"
# B4 keeps its existing text (Describe the nature of the business...)

$dpc.Range("A5").Value = "relationship"
$dpc.Range("B5").Value = "relationship"

$dpc.Range("A6").Value = "additional"
$dpc.Range("B6").Value = "Additional information "

$dpc.Hyperlinks.Add($dpc.Range("B3"), "mailto:ak@gmail.com") | Out-Null

$dpc.Activate()
$dpc.Range("D12").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Rename "demo" -> "rgth" and build the "Receive Gift Auto" form
# ---------------------------------------------------------------------
$rgth = $wb.Worksheets.Item("demo")
$rgth.Name = "rgth"

# NOTE: the engine snaps ColumnWidth to 1/6-character increments, so the
# inputs below are chosen so the stored width lands as close as possible
# to the original workbook's 30.140625 / 50.85546875 character widths.
$rgth.Columns.Item(1).ColumnWidth = 29.25062
$rgth.Columns.Item(2).ColumnWidth = 49.91747

$rgth.Range("A1").Value = "TestCases"
$rgth.Range("B1").Value = "Data"

$rgth.Range("A2").Value = "Request_Name"
$rgth.Range("B2").Value = "Receive Gift Auto"

$rgth.Range("A3").Value = "Full_Name"
$rgth.Range("B3").Value = "Akshay"

$rgth.Range("A4").Value = "Employer"
$rgth.Range("B4").Value = "Shreya"

$rgth.Range("A5").Value = "Offical_Position"
$rgth.Range("B5").Value = "QA"

$rgth.Range("A6").Value = "email_address"
$rgth.Range("B6").Value = "ak@gmail.com"
$rgth.Hyperlinks.Add($rgth.Range("B6"), "mailto:ak@gmail.com") | Out-Null

$rgth.Range("A7").Value = "Description"
$rgth.Range("B7").Value = "Details about the gift, travel or hospitality"

$rgth.Range("A8").Value = "Business_Purpose"
$rgth.Range("B8").Value = "Business Purpose or Rationale."

$rgth.Range("A9").Value = "Monetary_Value"
$rgth.Range("B9").Value = 500

$rgth.Activate()
$rgth.Range("A15").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Add the "pgth" sheet (Provide Gift Auto form) as the last sheet
# ---------------------------------------------------------------------
$pgth = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$pgth.Name = "pgth"

# See note above: chosen to land as close as possible to the target
# 24.5703125 / 28.85546875 character widths after the engine's rounding.
$pgth.Columns.Item(1).ColumnWidth = 23.58431
$pgth.Columns.Item(2).ColumnWidth = 27.91747

$pgth.Range("A1").Value = "TestCases"
$pgth.Range("B1").Value = "Data"

$pgth.Range("A2").Value = "RequestName"
$pgth.Range("B2").Value = "Provide Gift Auto"

$pgth.Range("A3").Value = "Describe the gift"
$pgth.Range("B3").Value = "hospitality in detail, including the business purpose or rationale"

$pgth.Range("A4").Value = "Additional Information"
$pgth.Range("B4").Value = "upload any additional files to support your approval request (optional)."

$pgth.Range("A5").Value = "Sort Description"
$pgth.Range("B5").Value = "Category Entertainment"

$pgth.Range("A6").Value = "Amount"
$pgth.Range("B6").Value = 500

$pgth.Range("A7").Value = "FullName"
$pgth.Range("B7").Value = "akshay"

$pgth.Range("A8").Value = "Employer"
$pgth.Range("B8").Value = "kapil"

$pgth.Range("A9").Value = "Title"
$pgth.Range("B9").Value = "Tester"

$pgth.Range("A10").Value = "Email"
$pgth.Range("B10").Value = "aagg@gmail.com"
$pgth.Hyperlinks.Add($pgth.Range("B10"), "mailto:aagg@gmail.com") | Out-Null

$pgth.Range("B6").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 4. Make "pgth" the active/selected sheet (activeTab = 3, tabSelected)
# ---------------------------------------------------------------------
$pgth.Activate()
$pgth.Range("D12").Select() | Out-Null
